# Fruta / hortaliza, semanal
# Insert a new weekly price-report group (3 rows: Especial/Primera/Segunda)
# at the top of the Mango data block (rows 337-339), pushing the rest of
# the existing rows down by 3 (old row 337 -> new row 340, ..., old row
# 377 -> new row 380). The new group uses the date 2021-09-10 (serial
# 44449) with the same price/origin pattern as the existing most-recent
# (2021-09-09 / Brasil) group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows above row 337, shifting rows 337:380 down to 340:380.
$ws.Range("A337:A339").EntireRow.Insert()

# Row 337 - Especial
$ws.Range("A337").Value = 8
$ws.Range("B337").Value = "Terminal La Palmera de La Serena"
$ws.Range("C337").Value = "Coquimbo"
$ws.Range("D337").Value = 44449
$ws.Range("E337").Value = 4
$ws.Range("F337").Value = "Fruta"
$ws.Range("G337").Value = 100108
$ws.Range("H337").Value = "Tropicales y subtropicales"
$ws.Range("I337").Value = 100108002
$ws.Range("J337").Value = "Mango"
$ws.Range("K337").Value = "Sin especificar"
$ws.Range("L337").Value = "Especial"
$ws.Range("M337").Value = 512
$ws.Range("N337").Value = 8000
$ws.Range("O337").Value = 8500
$ws.Range("P337").Value = 8250
$ws.Range("Q337").Value = '$/bandeja 4 kilos'
$ws.Range("R337").Value = 'Brasil'
$ws.Range("S337").Value = 2062
$ws.Range("T337").Value = 4

# Row 338 - Primera
$ws.Range("A338").Value = 8
$ws.Range("B338").Value = "Terminal La Palmera de La Serena"
$ws.Range("C338").Value = "Coquimbo"
$ws.Range("D338").Value = 44449
$ws.Range("E338").Value = 4
$ws.Range("F338").Value = "Fruta"
$ws.Range("G338").Value = 100108
$ws.Range("H338").Value = "Tropicales y subtropicales"
$ws.Range("I338").Value = 100108002
$ws.Range("J338").Value = "Mango"
$ws.Range("K338").Value = "Sin especificar"
$ws.Range("L338").Value = "Primera"
$ws.Range("M338").Value = 512
$ws.Range("N338").Value = 8000
$ws.Range("O338").Value = 8500
$ws.Range("P338").Value = 8250
$ws.Range("Q338").Value = '$/bandeja 4 kilos'
$ws.Range("R338").Value = 'Brasil'
$ws.Range("S338").Value = 2062
$ws.Range("T338").Value = 4

# Row 339 - Segunda
$ws.Range("A339").Value = 8
$ws.Range("B339").Value = "Terminal La Palmera de La Serena"
$ws.Range("C339").Value = "Coquimbo"
$ws.Range("D339").Value = 44449
$ws.Range("E339").Value = 4
$ws.Range("F339").Value = "Fruta"
$ws.Range("G339").Value = 100108
$ws.Range("H339").Value = "Tropicales y subtropicales"
$ws.Range("I339").Value = 100108002
$ws.Range("J339").Value = "Mango"
$ws.Range("K339").Value = "Sin especificar"
$ws.Range("L339").Value = "Segunda"
$ws.Range("M339").Value = 512
$ws.Range("N339").Value = 8000
$ws.Range("O339").Value = 8500
$ws.Range("P339").Value = 8250
$ws.Range("Q339").Value = '$/bandeja 4 kilos'
$ws.Range("R339").Value = 'Brasil'
$ws.Range("S339").Value = 2062
$ws.Range("T339").Value = 4
